$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1): update "想去人数" (col F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 7637
$ws1.Range("F11").Value = 643
$ws1.Range("F12").Value = 11

# Sheet "全部类型" (sheetId 4): same events appear one row lower, update col F too
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 7637
$ws4.Range("F12").Value = 643
$ws4.Range("F13").Value = 11
